# Applies the "assign3 params" tweak described in the commit:
#   "tweaks to assignment 3 and session 10 slides"
#
# On the "params" worksheet, remove the two rows describing the
# "True positive test result" (costTP, row 25) and "False positive
# test result" (costFP, row 26) parameters. All rows below shift up
# by two (old row 51 -> new row 49). Then refresh the AutoFilter
# range and the _FilterDatabase defined name so they reflect the new
# (smaller) table, and leave the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")
$ws.Activate()

# Row 25 = "True positive test result" / costTP / 91.8
# Row 26 = "False positive test result" / costFP / 91.8
# Deleting row 25 twice removes both (row 26 becomes row 25 after the
# first delete shifts everything up).
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(25).Delete()

# Recompute the AutoFilter over the now-smaller table (A1:E49 instead
# of A1:E51).
$ws.AutoFilterMode = $false
$ws.Range("A1:E49").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the
# shrunk filter range.
try {
    $wb.Names.Item("params!_FilterDatabase").RefersTo = "=params!`$A`$1:`$E`$49"
} catch {
}

# Move the selection to where the author left off editing.
$ws.Range("C20").Select()
